$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a number by Excel (loses exact text representation).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.585.25"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.562.58"
$ws.Range("E3").Value = "  +5.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "573.53"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("D6").Value = "150.57"
$ws.Range("E6").Value = "  +8.59%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "2.556.40"
$ws.Range("E9").Value = "  +4.90%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "28.24"
$ws.Range("E14").Value = "  +9.32%  "
$ws.Range("D15").Value = "3.018.48"
$ws.Range("E15").Value = "  +5.13%  "
$ws.Range("D16").Value = "63.527.50"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "2.578.88"
$ws.Range("E18").Value = "  +5.47%  "
$ws.Range("D19").Value = "11.71"
$ws.Range("E19").Value = "  +5.05%  "
$ws.Range("D20").Value = "342.72"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "6.95"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  +4.98%  "
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "1.48"
$ws.Range("E29").Value = "  +8.71%  "
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +13.82%  "
$ws.Range("D31").Value = "0.0₃0842"
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("D32").Value = "1.88"
$ws.Range("D33").Value = "177.58"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +8.73%  "
$ws.Range("D35").Value = "418.35"
$ws.Range("E35").Value = "  +11.62%  "
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "19.18"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "40.09"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").Value = "155.93"
$ws.Range("E43").Value = "  +6.29%  "
$ws.Range("D44").Value = "3.82"
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("D45").Value = "21.31"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").Value = "0.0969"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").Value = "18.87"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("D51").Value = "1.87"
$ws.Range("E51").Value = "  +9.15%  "
